# Apply odds updates to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("K2").Value = 2.05
$ws.Range("R2").Value = 1.57
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.5
$ws.Range("U2").Value = 2.25
$ws.Range("V2").Value = 1.57
$ws.Range("AC2").Value = 7
$ws.Range("AF2").Value = 81
$ws.Range("AT2").Value = 2.5
$ws.Range("AW2").Value = 7
$ws.Range("AX2").Value = 34

# Row 3
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.65

# Row 4
$ws.Range("K4").Value = 1.95
$ws.Range("L4").Value = 2.88
$ws.Range("R4").Value = 1.5

# Row 6
$ws.Range("V6").Value = 1.54

# Row 7
$ws.Range("V7").Value = 1.54

# Row 11
$ws.Range("M11").Value = 1.08
$ws.Range("O11").Value = 1.5
$ws.Range("P11").Value = 2.37

# Row 12
$ws.Range("G12").Value = 1.67
$ws.Range("I12").Value = 5.75
$ws.Range("M12").Value = 1.07
$ws.Range("O12").Value = 1.41
$ws.Range("P12").Value = 2.62
$ws.Range("U12").Value = 2.25
$ws.Range("V12").Value = 1.54
$ws.Range("AZ12").Value = 126

# Row 13
$ws.Range("M13").Value = 1.03
$ws.Range("O13").Value = 1.25
